$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row 3 (pushing old row 3 -> row 4, old row 4 -> row 5) for the
# "timestep 1" average row.
$ws.Rows.Item(3).Insert()

# Insert a new row 5 (pushing old row 4, now at row 5, -> row 6) for the
# "timestep 3" average row.
$ws.Rows.Item(5).Insert()

# Re-number the timestep labels in column A for every data row (0..4) and
# give the two new label cells the same formatting as the existing ones.
$ws.Cells.Item(2, 1).Copy()
$ws.Cells.Item(3, 1).PasteSpecial(-4122)
$ws.Cells.Item(5, 1).PasteSpecial(-4122)
$excel.CutCopyMode = 0

$ws.Cells.Item(2, 1).Value = 0
$ws.Cells.Item(3, 1).Value = 1
$ws.Cells.Item(4, 1).Value = 2
$ws.Cells.Item(5, 1).Value = 3
$ws.Cells.Item(6, 1).Value = 4

# Row 3: formulas that are the average of rows 2 and 4 (timestep 0 and 2).
$ws.Range("B3").Formula = "=(B2+B4)/2"
$ws.Range("C3:H3").Formula = "=(C2+C4)/2"

# Row 5: formulas that are the average of rows 4 and 6 (timestep 2 and 4).
$ws.Range("B5").Formula = "=(B4+B6)/2"
$ws.Range("C5:H5").Formula = "=(C4+C6)/2"

# The two newly-inserted data rows lose the special formatting on the
# "deaths"/"affected people" columns (E and G).
$ws.Cells.Item(3, 5).Style = "Normal"
$ws.Cells.Item(3, 7).Style = "Normal"
$ws.Cells.Item(5, 5).Style = "Normal"
$ws.Cells.Item(5, 7).Style = "Normal"

# Final selection landed a bit below the table, as seen in the saved file.
$ws.Range("H11").Select() | Out-Null

$wb.Save()
